# Update inputs mkt factor prem
$wb = $excel.ActiveWorkbook

# --- Sheet "fi": update duration/spread figures ---
$ws_fi = $wb.Worksheets.Item("fi")
$ws_fi.Range("B3").Value = 23.273078999999999
$ws_fi.Range("B4").Value = 12.668424
$ws_fi.Range("C4").Value = 0.01132061
$ws_fi.Range("D4").Value = 165.94323712832238
$ws_fi.Range("E4").Value = 50.524929419595544
$ws_fi.Range("F4").Value = 113.20609999999999
$ws_fi.Range("B5").Value = 16.484375

# --- Sheet "ret_assump": update TSY_MKT_RET return assumption ---
$ws_ret = $wb.Worksheets.Item("ret_assump")
$ws_ret.Range("B2").Value = 0.049428

# --- Sheet "mkt_factor_prem": update factor premium figures and highlight cells ---
$ws_mfp = $wb.Worksheets.Item("mkt_factor_prem")
$ws_mfp.Range("B8").Value = 0
$ws_mfp.Range("B9").Value = 0
$ws_mfp.Range("B10").Value = 0

$ws_mfp.Range("B5").Interior.Color = 65535
$ws_mfp.Range("B7").Interior.Color = 65535
$ws_mfp.Range("B8").Interior.Color = 65535
$ws_mfp.Range("B9").Interior.Color = 65535
$ws_mfp.Range("B10").Interior.Color = 65535

# --- Update window zoom levels and selections to match author's final view state ---
$ws_fi.Select()
$excel.ActiveWindow.Zoom = 130
$ws_fi.Range("B25").Select()

$ws_ret.Select()
$ws_ret.Range("B2").Select()

$ws_mfp.Select()
$excel.ActiveWindow.Zoom = 115
$ws_mfp.Range("D9").Select()
